# BL Audit Form update: "this added last date 21-11-24"
#
# - Report header date (B1) moves from 17.11.2024 -> 21.11.2024
# - Stock purchase quantities for a few product rows are updated
# - Two "liability" entries (E22/E23) are corrected
# - A new Sim commission entry "E-life Comm" / 70000 is recorded (row 31)
# - The bank-guarantee credit amount and its note move from
#   100000 / "17.11.2024 payment" to 400000 / "24.11.2024 payment"
# All dependent SUM()/formula cells recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -------------------------------------------------------------
$ws.Range("B1").Value = "21.11.2024"

# --- Stock table (rows 9-21) ---------------------------------------------
$ws.Range("C9").Value  = 393510   # was 308382
$ws.Range("C10").Value = 1000     # was blank
$ws.Range("C12").Value = 18000    # was 12620
$ws.Range("C16").Value = 75       # was 18

# --- Receivable / liability rows -----------------------------------------
$ws.Range("E22").Value = 32554    # was 14519
$ws.Range("E23").Value = 12212    # was 23022

# --- New sim commission entry (row 31) ------------------------------------
$ws.Range("D31").Value = "E-life Comm"
$ws.Range("E31").Value = 70000

# --- Bank guarantee credit row (row 34) -----------------------------------
$ws.Range("E34").Value = 400000
$ws.Range("F34").Value = "24.11.2024 payment "

# --- Update the on-screen selection to match the saved view --------------
$ws.Range("F35").Select()
